$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values (after edit)
$ws.Range("D2").Value = 44350
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 167

# Row 4 values (after edit)
$ws.Range("D4").Value = 44273
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 233
